# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) across all
# 8 Disciple-of-the-Hand sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with
# refreshed Universalis market-board figures from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 156.46153
$ws.Range("I9").Value = 121.77778
$ws.Range("J9").Value = 234.5
$ws.Range("K9").Value = 121.77778
$ws.Range("L9").Value = 234.5
$ws.Range("M9").Value = 47.22221999999999
$ws.Range("N9").Value = -572.5
$ws.Range("H40").Value = 2025
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H43").Value = 37533.332
$ws.Range("I43").Value = 49999.668
$ws.Range("J43").Value = 31300.166
$ws.Range("K43").Value = 49999.668
$ws.Range("L43").Value = 31300.166
$ws.Range("M43").Value = -49930.668
$ws.Range("N43").Value = -31438.166
$ws.Range("H70").Value = 9672.25
$ws.Range("J70").Value = 10625.571
$ws.Range("L70").Value = 31876.713
$ws.Range("N70").Value = -32416.713
$ws.Range("H73").Value = 9672.25
$ws.Range("J73").Value = 10625.571
$ws.Range("L73").Value = 31876.713
$ws.Range("N73").Value = -33748.713
$ws.Range("H94").Value = 581.5714
$ws.Range("I94").Value = 581.5714
$ws.Range("K94").Value = 581.5714
$ws.Range("M94").Value = -130.5714
$ws.Range("H132").Value = 13979.5
$ws.Range("I132").Value = 14331.235
$ws.Range("K132").Value = 42993.705
$ws.Range("M132").Value = -40463.705
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 693.125
$ws.Range("I2").Value = 506.42856
$ws.Range("K2").Value = 506.42856
$ws.Range("M2").Value = -393.42856
$ws.Range("H6").Value = 751500
$ws.Range("I6").Value = 751500
$ws.Range("K6").Value = 751500
$ws.Range("M6").Value = -751327
$ws.Range("H32").Value = 125308.58
$ws.Range("I32").Value = 130893.86
$ws.Range("J32").Value = 14999.25
$ws.Range("K32").Value = 130893.86
$ws.Range("L32").Value = 14999.25
$ws.Range("M32").Value = -130606.86
$ws.Range("N32").Value = -15573.25
$ws.Range("H45").Value = 1971.2142
$ws.Range("J45").Value = 1966.5
$ws.Range("L45").Value = 1966.5
$ws.Range("N45").Value = -2720.5
$ws.Range("H61").Value = 7074.067
$ws.Range("I61").Value = 8871.6
$ws.Range("J61").Value = 3479
$ws.Range("K61").Value = 8871.6
$ws.Range("L61").Value = 3479
$ws.Range("M61").Value = -8659.6
$ws.Range("N61").Value = -3903
$ws.Range("H74").Value = 12617.85
$ws.Range("I74").Value = 1530.5834
$ws.Range("J74").Value = 29248.75
$ws.Range("K74").Value = 1530.5834
$ws.Range("L74").Value = 29248.75
$ws.Range("M74").Value = -656.5834
$ws.Range("N74").Value = -30996.75
$ws.Range("H77").Value = 12617.85
$ws.Range("I77").Value = 1530.5834
$ws.Range("J77").Value = 29248.75
$ws.Range("K77").Value = 7652.916999999999
$ws.Range("L77").Value = 146243.75
$ws.Range("M77").Value = -3284.916999999999
$ws.Range("N77").Value = -154979.75
$ws.Range("H110").Value = 1428.8889
$ws.Range("I110").Value = 1143.3334
$ws.Range("K110").Value = 1143.3334
$ws.Range("M110").Value = 901.6666
$ws.Range("H116").Value = 693.125
$ws.Range("I116").Value = 506.42856
$ws.Range("K116").Value = 506.42856
$ws.Range("M116").Value = 1787.57144
$ws.Range("H122").Value = 1519.4615
$ws.Range("I122").Value = 1025.4
$ws.Range("J122").Value = 3166.3333
$ws.Range("K122").Value = 3076.2
$ws.Range("L122").Value = 9498.999899999999
$ws.Range("M122").Value = -626.2000000000003
$ws.Range("N122").Value = -14398.9999
$ws.Range("H132").Value = 1391200.5
$ws.Range("I132").Value = 1668820.6
$ws.Range("K132").Value = 5006461.800000001
$ws.Range("M132").Value = -5003931.800000001
$ws.Range("H134").Value = 89999
$ws.Range("J134").Value = 89999
$ws.Range("L134").Value = 89999
$ws.Range("N134").Value = -100139
$ws.Range("H136").Value = 7074.067
$ws.Range("I136").Value = 8871.6
$ws.Range("J136").Value = 3479
$ws.Range("K136").Value = 26614.8
$ws.Range("L136").Value = 10437
$ws.Range("M136").Value = -24064.8
$ws.Range("N136").Value = -15537
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 693.125
$ws.Range("I3").Value = 506.42856
$ws.Range("K3").Value = 506.42856
$ws.Range("M3").Value = -392.42856
$ws.Range("H134").Value = 9979.954
$ws.Range("I134").Value = 6093.9473
$ws.Range("J134").Value = 34591.332
$ws.Range("K134").Value = 18281.8419
$ws.Range("L134").Value = 103773.996
$ws.Range("M134").Value = -15746.8419
$ws.Range("N134").Value = -108843.996
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1417.2941
$ws.Range("I16").Value = 1505.7858
$ws.Range("J16").Value = 1004.3333
$ws.Range("K16").Value = 1505.7858
$ws.Range("L16").Value = 1004.3333
$ws.Range("M16").Value = -1218.7858
$ws.Range("N16").Value = -1578.3333
$ws.Range("H22").Value = 1488.8948
$ws.Range("I22").Value = 608.0833
$ws.Range("K22").Value = 608.0833
$ws.Range("M22").Value = -258.0833
$ws.Range("H31").Value = 3414.6511
$ws.Range("I31").Value = 4076.1667
$ws.Range("J31").Value = 2579.0527
$ws.Range("K31").Value = 4076.1667
$ws.Range("L31").Value = 2579.0527
$ws.Range("M31").Value = -3781.1667
$ws.Range("N31").Value = -3169.0527
$ws.Range("H34").Value = 3414.6511
$ws.Range("I34").Value = 4076.1667
$ws.Range("J34").Value = 2579.0527
$ws.Range("K34").Value = 4076.1667
$ws.Range("L34").Value = 2579.0527
$ws.Range("M34").Value = -3874.1667
$ws.Range("N34").Value = -2983.0527
$ws.Range("H62").Value = 3732.6365
$ws.Range("I62").Value = 2749.75
$ws.Range("K62").Value = 2749.75
$ws.Range("M62").Value = -2125.75
$ws.Range("H65").Value = 3732.6365
$ws.Range("I65").Value = 2749.75
$ws.Range("K65").Value = 13748.75
$ws.Range("M65").Value = -10628.75
$ws.Range("H113").Value = 1417.2941
$ws.Range("I113").Value = 1505.7858
$ws.Range("J113").Value = 1004.3333
$ws.Range("K113").Value = 1505.7858
$ws.Range("L113").Value = 1004.3333
$ws.Range("M113").Value = 664.2141999999999
$ws.Range("N113").Value = -5344.3333
$ws.Range("H132").Value = 3594
$ws.Range("I132").Value = 3594
$ws.Range("K132").Value = 10782
$ws.Range("M132").Value = -8252
$ws.Range("H134").Value = 3278.55
$ws.Range("I134").Value = 2670.611
$ws.Range("J134").Value = 8750
$ws.Range("K134").Value = 8011.833
$ws.Range("L134").Value = 26250
$ws.Range("M134").Value = -5476.833
$ws.Range("N134").Value = -31320
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 579.1429000000001
$ws.Range("J17").Value = 1230
$ws.Range("L17").Value = 3690
$ws.Range("N17").Value = -4028
$ws.Range("H48").Value = 699.3333
$ws.Range("I48").Value = 851.5
$ws.Range("K48").Value = 2554.5
$ws.Range("M48").Value = -2304.5
$ws.Range("H50").Value = 219467.95
$ws.Range("J50").Value = 835522.3
$ws.Range("L50").Value = 2506566.9
$ws.Range("N50").Value = -2507528.9
$ws.Range("H53").Value = 219467.95
$ws.Range("J53").Value = 835522.3
$ws.Range("L53").Value = 2506566.9
$ws.Range("N53").Value = -2507528.9
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3106.6155
$ws.Range("I113").Value = 2340
$ws.Range("K113").Value = 2340
$ws.Range("M113").Value = -170
$ws.Range("H126").Value = 6500.0625
$ws.Range("I126").Value = 6928.7144
$ws.Range("K126").Value = 20786.1432
$ws.Range("M126").Value = -18316.1432
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4440.75
$ws.Range("I7").Value = 4378.231
$ws.Range("K7").Value = 4378.231
$ws.Range("M7").Value = -4266.231
$ws.Range("H40").Value = 6302.4287
$ws.Range("I40").Value = 8055.5
$ws.Range("J40").Value = 3965
$ws.Range("K40").Value = 8055.5
$ws.Range("L40").Value = 3965
$ws.Range("M40").Value = -7919.5
$ws.Range("N40").Value = -4237
$ws.Range("H46").Value = 2729.85
$ws.Range("I46").Value = 399.75
$ws.Range("K46").Value = 399.75
$ws.Range("M46").Value = -211.75
$ws.Range("H61").Value = 11442.277
$ws.Range("I61").Value = 9821.352999999999
$ws.Range("K61").Value = 9821.352999999999
$ws.Range("M61").Value = -9619.352999999999
$ws.Range("H113").Value = 11442.277
$ws.Range("I113").Value = 9821.352999999999
$ws.Range("K113").Value = 9821.352999999999
$ws.Range("M113").Value = -7651.352999999999
$ws.Range("H126").Value = 4440.75
$ws.Range("I126").Value = 4378.231
$ws.Range("K126").Value = 13134.693
$ws.Range("M126").Value = -10664.693
$ws.Range("H132").Value = 3863785.2
$ws.Range("I132").Value = 5579321
$ws.Range("K132").Value = 16737963
$ws.Range("M132").Value = -16735433
$ws.Range("H136").Value = 8261.565000000001
$ws.Range("I136").Value = 4371.5835
$ws.Range("K136").Value = 13114.7505
$ws.Range("M136").Value = -10564.7505
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2523.926
$ws.Range("I132").Value = 2389.7273
$ws.Range("J132").Value = 3114.4
$ws.Range("K132").Value = 7169.1819
$ws.Range("L132").Value = 9343.200000000001
$ws.Range("M132").Value = -4639.1819
$ws.Range("N132").Value = -14403.2
$ws.Range("H136").Value = 893.5625
$ws.Range("I136").Value = 919.9091
$ws.Range("J136").Value = 835.6
$ws.Range("K136").Value = 2759.7273
$ws.Range("L136").Value = 2506.8
$ws.Range("M136").Value = -209.7273
$ws.Range("N136").Value = -7606.8
